$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update candidate row (row 2) with new user details
$ws.Range("A2").Value = 'test825'
$ws.Range("B2").Value = 23071033
$ws.Range("C2").Value = 'narendra91'
$ws.Range("D2").Value = 's5U%8$fB'
$ws.Range("E2").Value = 'MR'
$ws.Range("F2").Value = 'Narendra'
$ws.Range("G2").Value = 'Modi'
